# Re-order the batting lineup so each row's stats "follow" the player.
# Qaiser Patel moves up from row 15 to row 4 (right after Joe Edwards),
# and the players previously in rows 4, 8, 9, 10, 12, 14 each shift down
# one lineup slot (landing in rows 8, 9, 10, 12, 14, 15 respectively).
# Every other player/row (2,3,5,6,7,11,13,16,17,18) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 -> becomes Qaiser Patel (previously row 15's player/stats)
$ws.Range("A4").Value = "Qaiser Patel"
$ws.Range("B4").Value = 12
$ws.Range("C4").Value = 9
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 1
$ws.Range("N4").Value = 0.778
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 0.75
$ws.Range("Q4").Value = 1.75

# Row 8 -> becomes Rich Squitieri (previously row 4's player/stats)
$ws.Range("A8").Value = "Rich Squitieri"
$ws.Range("B8").Value = 21
$ws.Range("C8").Value = 21
$ws.Range("D8").Value = 4
$ws.Range("E8").Value = 14
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 11
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").Value = 0.667
$ws.Range("O8").Value = 0.857
$ws.Range("P8").Value = 0.667
$ws.Range("Q8").Value = 1.524

# Row 9 -> becomes Nick Mirman (previously row 8's player/stats)
$ws.Range("A9").Value = "Nick Mirman"
$ws.Range("B9").Value = 9
$ws.Range("C9").Value = 8
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 4
$ws.Range("J9").Value = 1
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").Value = 0.625
$ws.Range("O9").Value = 0.625
$ws.Range("P9").Value = 0.667
$ws.Range("Q9").Value = 1.292

# Row 10 -> becomes Nick Hanten (previously row 9's player/stats)
$ws.Range("A10").Value = "Nick Hanten"
$ws.Range("B10").Value = 18
$ws.Range("C10").Value = 18
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 4
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").Value = 0.5
$ws.Range("O10").Value = 0.5
$ws.Range("P10").Value = 0.5
$ws.Range("Q10").Value = 1

# Row 12 -> becomes Scott Richardson (previously row 10's player/stats)
$ws.Range("A12").Value = "Scott Richardson"
$ws.Range("B12").Value = 17
$ws.Range("C12").Value = 17
$ws.Range("D12").Value = 6
$ws.Range("E12").Value = 9
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").Value = 0.529
$ws.Range("O12").Value = 0.765
$ws.Range("P12").Value = 0.529
$ws.Range("Q12").Value = 1.294

# Row 14 -> becomes Gordon Walker (previously row 12's player/stats)
$ws.Range("A14").Value = "Gordon Walker"
$ws.Range("B14").Value = 19
$ws.Range("C14").Value = 17
$ws.Range("D14").Value = 10
$ws.Range("E14").Value = 9
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 5
$ws.Range("J14").Value = 2
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = 0.529
$ws.Range("O14").Value = 0.529
$ws.Range("P14").Value = 0.579
$ws.Range("Q14").Value = 1.108

# Row 15 -> becomes Andrew Scott (previously row 14's player/stats)
$ws.Range("A15").Value = "Andrew Scott"
$ws.Range("B15").Value = 7
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = 0.571
$ws.Range("O15").Value = 0.571
$ws.Range("P15").Value = 0.571
$ws.Range("Q15").Value = 1.143

